$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 141.70809
$ws.Cells.Item(2, 8).Value = 425.12427
$ws.Cells.Item(2, 9).Value = 0.4270657810795758
$ws.Cells.Item(2, 10).Value = 0.4270657810795759
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 127.984071
$ws.Cells.Item(2, 14).Value = 383.952213
$ws.Cells.Item(2, 15).Value = 0.4594739502473105
$ws.Cells.Item(2, 16).Value = 0.4594739502473105
$ws.Cells.Item(2, 17).Value = 18136.37825183439
$ws.Cells.Item(2, 18).Value = 163227.4042665095
$ws.Cells.Item(2, 19).Value = 0.1962256014480858
$ws.Cells.Item(2, 20).Value = 0.1962256014480859

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 141.70809
$ws.Cells.Item(3, 8).Value = 425.12427
$ws.Cells.Item(3, 9).Value = 0.4270657810795758
$ws.Cells.Item(3, 10).Value = 0.4270657810795759
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 60.45343933333334
$ws.Cells.Item(3, 14).Value = 181.360318
$ws.Cells.Item(3, 15).Value = 0.2170331070069088
$ws.Cells.Item(3, 16).Value = 0.2170331070069088
$ws.Cells.Item(3, 17).Value = 8566.74142185754
$ws.Cells.Item(3, 18).Value = 77100.67279671787
$ws.Cells.Item(3, 19).Value = 0.09268741336403265
$ws.Cells.Item(3, 20).Value = 0.09268741336403268

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 141.70809
$ws.Cells.Item(4, 8).Value = 425.12427
$ws.Cells.Item(4, 9).Value = 0.4270657810795758
$ws.Cells.Item(4, 10).Value = 0.4270657810795759
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 59.37981533333334
$ws.Cells.Item(4, 14).Value = 178.139446
$ws.Cells.Item(4, 15).Value = 0.2131787034353868
$ws.Cells.Item(4, 16).Value = 0.2131787034353868
$ws.Cells.Item(4, 17).Value = 8414.600215439381
$ws.Cells.Item(4, 18).Value = 75731.40193895443
$ws.Cells.Item(4, 19).Value = 0.09104132949216473
$ws.Cells.Item(4, 20).Value = 0.09104132949216474

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 141.70809
$ws.Cells.Item(5, 8).Value = 425.12427
$ws.Cells.Item(5, 9).Value = 0.4270657810795758
$ws.Cells.Item(5, 10).Value = 0.4270657810795759
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 28.25780433333334
$ws.Cells.Item(5, 14).Value = 84.77341300000001
$ws.Cells.Item(5, 15).Value = 0.1014479761497213
$ws.Cells.Item(5, 16).Value = 0.1014479761497213
$ws.Cells.Item(5, 17).Value = 4004.359479670391
$ws.Cells.Item(5, 18).Value = 36039.23531703351
$ws.Cells.Item(5, 19).Value = 0.04332495917332291
$ws.Cells.Item(5, 20).Value = 0.04332495917332291

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 141.70809
$ws.Cells.Item(6, 8).Value = 425.12427
$ws.Cells.Item(6, 9).Value = 0.4270657810795758
$ws.Cells.Item(6, 10).Value = 0.4270657810795759
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 2.469651333333333
$ws.Cells.Item(6, 14).Value = 7.408954
$ws.Cells.Item(6, 15).Value = 0.008866263160672582
$ws.Cells.Item(6, 16).Value = 0.008866263160672582
$ws.Cells.Item(6, 17).Value = 349.96957341262
$ws.Cells.Item(6, 18).Value = 3149.72616071358
$ws.Cells.Item(6, 19).Value = 0.003786477601969705
$ws.Cells.Item(6, 20).Value = 0.003786477601969706

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 171.783722
$ws.Cells.Item(7, 8).Value = 515.3511659999999
$ws.Cells.Item(7, 9).Value = 0.5177047366363254
$ws.Cells.Item(7, 10).Value = 0.5177047366363255
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 127.984071
$ws.Cells.Item(7, 14).Value = 383.952213
$ws.Cells.Item(7, 15).Value = 0.4594739502473105
$ws.Cells.Item(7, 16).Value = 0.4594739502473105
$ws.Cells.Item(7, 17).Value = 21985.58007309226
$ws.Cells.Item(7, 18).Value = 197870.2206578303
$ws.Cells.Item(7, 19).Value = 0.237871840404036
$ws.Cells.Item(7, 20).Value = 0.2378718404040361

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 171.783722
$ws.Cells.Item(8, 8).Value = 515.3511659999999
$ws.Cells.Item(8, 9).Value = 0.5177047366363254
$ws.Cells.Item(8, 10).Value = 0.5177047366363255
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 60.45343933333334
$ws.Cells.Item(8, 14).Value = 181.360318
$ws.Cells.Item(8, 15).Value = 0.2170331070069088
$ws.Cells.Item(8, 16).Value = 0.2170331070069088
$ws.Cells.Item(8, 17).Value = 10384.9168163812
$ws.Cells.Item(8, 18).Value = 93464.25134743078
$ws.Cells.Item(8, 19).Value = 0.1123590675043751
$ws.Cells.Item(8, 20).Value = 0.1123590675043752

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 171.783722
$ws.Cells.Item(9, 8).Value = 515.3511659999999
$ws.Cells.Item(9, 9).Value = 0.5177047366363254
$ws.Cells.Item(9, 10).Value = 0.5177047366363255
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 59.37981533333334
$ws.Cells.Item(9, 14).Value = 178.139446
$ws.Cells.Item(9, 15).Value = 0.2131787034353868
$ws.Cells.Item(9, 16).Value = 0.2131787034353868
$ws.Cells.Item(9, 17).Value = 10200.48568963267
$ws.Cells.Item(9, 18).Value = 91804.37120669404
$ws.Cells.Item(9, 19).Value = 0.1103636245184903
$ws.Cells.Item(9, 20).Value = 0.1103636245184903

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 171.783722
$ws.Cells.Item(10, 8).Value = 515.3511659999999
$ws.Cells.Item(10, 9).Value = 0.5177047366363254
$ws.Cells.Item(10, 10).Value = 0.5177047366363255
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 28.25780433333334
$ws.Cells.Item(10, 14).Value = 84.77341300000001
$ws.Cells.Item(10, 15).Value = 0.1014479761497213
$ws.Cells.Item(10, 16).Value = 0.1014479761497213
$ws.Cells.Item(10, 17).Value = 4854.230803927729
$ws.Cells.Item(10, 18).Value = 43688.07723534956
$ws.Cells.Item(10, 19).Value = 0.05252009777487969
$ws.Cells.Item(10, 20).Value = 0.0525200977748797

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 171.783722
$ws.Cells.Item(11, 8).Value = 515.3511659999999
$ws.Cells.Item(11, 9).Value = 0.5177047366363254
$ws.Cells.Item(11, 10).Value = 0.5177047366363255
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 2.469651333333333
$ws.Cells.Item(11, 14).Value = 7.408954
$ws.Cells.Item(11, 15).Value = 0.008866263160672582
$ws.Cells.Item(11, 16).Value = 0.008866263160672582
$ws.Cells.Item(11, 17).Value = 424.2458980822626
$ws.Cells.Item(11, 18).Value = 3818.213082740363
$ws.Cells.Item(11, 19).Value = 0.004590106434544354
$ws.Cells.Item(11, 20).Value = 0.004590106434544355

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 18.32613333333333
$ws.Cells.Item(12, 8).Value = 54.9784
$ws.Cells.Item(12, 9).Value = 0.05522948228409861
$ws.Cells.Item(12, 10).Value = 0.05522948228409861
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 127.984071
$ws.Cells.Item(12, 14).Value = 383.952213
$ws.Cells.Item(12, 15).Value = 0.4594739502473105
$ws.Cells.Item(12, 16).Value = 0.4594739502473105
$ws.Cells.Item(12, 17).Value = 2345.4531496888
$ws.Cells.Item(12, 18).Value = 21109.0783471992
$ws.Cells.Item(12, 19).Value = 0.02537650839518864
$ws.Cells.Item(12, 20).Value = 0.02537650839518864

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 18.32613333333333
$ws.Cells.Item(13, 8).Value = 54.9784
$ws.Cells.Item(13, 9).Value = 0.05522948228409861
$ws.Cells.Item(13, 10).Value = 0.05522948228409861
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 60.45343933333334
$ws.Cells.Item(13, 14).Value = 181.360318
$ws.Cells.Item(13, 15).Value = 0.2170331070069088
$ws.Cells.Item(13, 16).Value = 0.2170331070069088
$ws.Cells.Item(13, 17).Value = 1107.877789681245
$ws.Cells.Item(13, 18).Value = 9970.900107131201
$ws.Cells.Item(13, 19).Value = 0.01198662613850095
$ws.Cells.Item(13, 20).Value = 0.01198662613850095

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 18.32613333333333
$ws.Cells.Item(14, 8).Value = 54.9784
$ws.Cells.Item(14, 9).Value = 0.05522948228409861
$ws.Cells.Item(14, 10).Value = 0.05522948228409861
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 59.37981533333334
$ws.Cells.Item(14, 14).Value = 178.139446
$ws.Cells.Item(14, 15).Value = 0.2131787034353868
$ws.Cells.Item(14, 16).Value = 0.2131787034353868
$ws.Cells.Item(14, 17).Value = 1088.202413107378
$ws.Cells.Item(14, 18).Value = 9793.821717966401
$ws.Cells.Item(14, 19).Value = 0.01177374942473181
$ws.Cells.Item(14, 20).Value = 0.01177374942473181

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 18.32613333333333
$ws.Cells.Item(15, 8).Value = 54.9784
$ws.Cells.Item(15, 9).Value = 0.05522948228409861
$ws.Cells.Item(15, 10).Value = 0.05522948228409861
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 28.25780433333334
$ws.Cells.Item(15, 14).Value = 84.77341300000001
$ws.Cells.Item(15, 15).Value = 0.1014479761497213
$ws.Cells.Item(15, 16).Value = 0.1014479761497213
$ws.Cells.Item(15, 17).Value = 517.8562899199112
$ws.Cells.Item(15, 18).Value = 4660.7066092792
$ws.Cells.Item(15, 19).Value = 0.005602919201518691
$ws.Cells.Item(15, 20).Value = 0.005602919201518691

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 18.32613333333333
$ws.Cells.Item(16, 8).Value = 54.9784
$ws.Cells.Item(16, 9).Value = 0.05522948228409861
$ws.Cells.Item(16, 10).Value = 0.05522948228409861
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 2.469651333333333
$ws.Cells.Item(16, 14).Value = 7.408954
$ws.Cells.Item(16, 15).Value = 0.008866263160672582
$ws.Cells.Item(16, 16).Value = 0.008866263160672582
$ws.Cells.Item(16, 17).Value = 45.25915962151111
$ws.Cells.Item(16, 18).Value = 407.3324365936
$ws.Cells.Item(16, 19).Value = 0.0004896791241585225
$ws.Cells.Item(16, 20).Value = 0.0004896791241585225
